$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook level: rename the second sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "procfs_sam"

# ---------------------------------------------------------------------
# Sheet 1 ("procfs"): only the selected cell and row 4's height change;
# none of its cell contents moved. Select this before touching sheet 2
# so sheet 2 ends up as the active tab (matching the source file).
# ---------------------------------------------------------------------
$ws1.Rows.Item(4).RowHeight = 240
$ws1.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2 ("procfs_sam"): content got reshuffled/renamed (ylmtest -> test_proc),
# gained an extra header/divider line plus a couple of new statements, and
# moved onto extra columns (C/D) for the nested code blocks. Easiest and
# most reliable way to reproduce this faithfully is to wipe the sheet and
# re-author every used cell from scratch in its final layout.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Cells.Clear()

$ws2.Range("B2").Value = "################################ procfs ################################"
$ws2.Range("B3").Value = "1. 创建proc文件夹"
$ws2.Range("C4").Value = "struct proc_dir_entry *test_proc_dir = proc_mkdir(""test_proc_dir"", NULL);"

$ws2.Range("B6").Value = "2. 创建proc文件"
$ws2.Range("C7").Value = "static int test_proc_read(struct seq_file * m, void* vp);"
$ws2.Range("D8").Value = "seq_puts()"
$ws2.Range("D9").Value = "seq_printf()"

$ws2.Range("C11").Value = "static ssize_t test_proc_write(struct file *file, const char *buffer, size_t count, loff_t *loff);"
$ws2.Range("D12").Value = "copy_from_user()"

$ws2.Range("C14").Value = "static int test_proc_open(struct inode *inode, struct file *file);"
$ws2.Range("D15").Value = "return single_open(file, &test_proc_read, NULL);"

$ws2.Range("C17").Value = "static struct file_operations test_proc_fops = {"
$ws2.Range("C18").Value = "    .open = test_proc_open,"
$ws2.Range("C19").Value = "    .read = seq_read,"
$ws2.Range("C20").Value = "    .write = test_proc_write,"
$ws2.Range("C21").Value = "};"

$ws2.Range("C23").Value = "struct proc_dir_entry *test_proc_entry;"
$ws2.Range("C24").Value = "test_proc_entry = proc_create_data(""test_proc"", 0666, test_proc_dir, &test_proc_fops, NULL);"

$ws2.Range("B26").Value = "3. 移除proc文件夹或文件"
$ws2.Range("C27").Value = "remove_proc_entry(""test_proc"", ""test_proc_dir"");"
$ws2.Range("C28").Value = "remove_proc_entry(""test_proc_dir"", NULL);"

$win = $excel.ActiveWindow
$win.DisplayGridlines = $false
$ws2.Range("AB23").Select() | Out-Null
